$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force-set a cell value as text even if it looks numeric,
# while preserving the cells original style (no NumberFormat/style drift).

$ws.Range("D2").Value = '26.891.81'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '1.816.32'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -0.02%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.21'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -0.48%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +0.01%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4689'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  +1.41%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3694'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  -1.66%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07376'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  -0.53%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8707'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +1.00%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.44'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = '1.821.58'
$ws.Range("E12").Value = '  +0.46%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.368'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -0.36%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.38'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +0.50%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07085'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -0.06%  '
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.504'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  -2.26%  '
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("E20").Value = '  -1.02%  '
$ws.Range("D21").Value = '26.911.94'
$ws.Range("E21").Value = '  -0.76%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.341'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("E23").Value = '  -2.86%  '
$ws.Range("D24").Value = '2.038.91'
$ws.Range("E24").Value = '  -0.44%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.895'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  -1.44%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.71'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  +0.08%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.199'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -0.61%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.42'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -0.51%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.305'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +0.48%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.52'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  -1.25%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08930'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +0.03%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7673'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -0.69%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.161'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -0.79%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.479'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  -1.13%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.923'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +1.28%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.001'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  +0.02%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.100'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -2.28%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01958'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -0.11%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05262'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +0.63%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.943'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.278'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5331'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +0.92%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.364'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -0.44%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1661'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -1.01%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.454'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -1.81%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4952'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  -1.78%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.43'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +0.89%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.672'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +0.03%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.83'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -2.12%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06286'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -0.51%  '
